$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.533.20'
$ws.Range("E2").Value = '  -1.89%  '
$ws.Range("D3").Value = '2.628.87'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '577.40'
$ws.Range("E5").Value = '  -3.67%  '
$ws.Range("D6").Value = '156.12'
$ws.Range("D7").Value = '0.646'
$ws.Range("E7").Value = '  +5.27%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  -5.23%  '
$ws.Range("D10").Value = '5.80'
$ws.Range("E10").Value = '  -0.61%  '
$ws.Range("D11").Value = '0.389'
$ws.Range("E11").Value = '  -2.65%  '
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '28.53'
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("E14").Value = '  -6.51%  '
$ws.Range("D15").Value = '3.100.40'
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("D16").Value = '64.338.64'
$ws.Range("E16").Value = '  -2.04%  '
$ws.Range("D17").Value = '2.619.41'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = '12.23'
$ws.Range("E18").Value = '  -3.85%  '
$ws.Range("D19").Value = '4.69'
$ws.Range("E19").Value = '  -2.36%  '
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").Value = '  -1.93%  '
$ws.Range("D21").Value = '345.77'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").Value = '67.89'
$ws.Range("E23").Value = '  -2.38%  '
$ws.Range("E24").Value = '  -1.98%  '
$ws.Range("E25").Value = '  +3.24%  '
$ws.Range("D26").Value = '9.35'
$ws.Range("E26").Value = '  -4.08%  '
$ws.Range("D27").Value = '1.55'
$ws.Range("E27").Value = '  -3.48%  '
$ws.Range("D28").Value = '551.26'
$ws.Range("E28").Value = '  +3.73%  '
$ws.Range("D29").Value = '0.161'
$ws.Range("E29").Value = '  -2.61%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.05%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '7.92'
$ws.Range("E31").Value = '  -1.42%  '
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -3.56%  '
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("D34").Value = '6.42'
$ws.Range("E34").Value = '  -1.00%  '
$ws.Range("D35").Value = '5.33'
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("E36").Value = '  -2.54%  '
$ws.Range("D37").Value = '20.01'
$ws.Range("E37").Value = '  -2.80%  '
$ws.Range("E38").Value = '  -0.05%  '
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("D40").Value = '152.62'
$ws.Range("E40").Value = '  -3.26%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("D42").Value = '2.47'
$ws.Range("E42").Value = '  +4.27%  '
$ws.Range("D43").Value = '158.22'
$ws.Range("E43").Value = '  -3.11%  '
$ws.Range("E45").Value = '  -2.03%  '
$ws.Range("D46").Value = '22.80'
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("E48").Value = '  +2.24%  '
$ws.Range("E49").Value = '  -3.21%  '
$ws.Range("D50").Value = '19.15'
$ws.Range("E50").Value = '  -4.57%  '
$ws.Range("D51").Value = '0.0₆0239'
$ws.Range("E51").Value = '  -6.81%  '
